# Inseriti dati mancanti nel dataset delle nascite
# Fills in previously-missing ("..") values in the births dataset with the
# real figures, keeping the cells as text (matching the rest of the sheet,
# which stores all numeric-looking values as text/inline strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "M2"  = "298.5"
    "M5"  = "373.7"
    "M6"  = "229.1"
    "M21" = "818.5"
    "M22" = "289"
    "M26" = "1880"
    "M38" = "677.2"
    "B40" = "734.5"
    "C40" = "739.3"
    "D40" = "743.8"
    "E40" = "747.4"
    "F40" = "749.6"
    "G40" = "750.3"
    "L40" = "636"
    "M40" = "629.4"
    "K41" = "2890"
    "M41" = "2760"
    "B45" = "26600"
    "C45" = "26340"
    "D45" = "26030"
    "E45" = "25740"
    "F45" = "24900"
    "G45" = "24830"
    "L45" = "23140"
    "M45" = "23110"
    "B46" = "617.6"
    "C46" = "613.2"
    "D46" = "602.8"
    "E46" = "592.9"
    "F46" = "584.7"
    "G46" = "577.6"
    "G48" = "1940"
    "H48" = "1890"
    "I48" = "1690"
    "J48" = "1610"
    "K48" = "1490"
    "L48" = "1440"
    "M48" = "1400"
    "M49" = "1180"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
}
